$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 613.7143
$ws.Range("I31").Value = 324
$ws.Range("K31").Value = 972
$ws.Range("M31").Value = -742

$ws.Range("H40").Value = 2084.5518
$ws.Range("I40").Value = 2857
$ws.Range("J40").Value = 1539.2941
$ws.Range("K40").Value = 2857
$ws.Range("L40").Value = 1539.2941
$ws.Range("M40").Value = -2682
$ws.Range("N40").Value = -1889.2941

$ws.Range("H128").Value = 39866.668
$ws.Range("J128").Value = 39866.668
$ws.Range("L128").Value = 39866.668
$ws.Range("N128").Value = -49826.668

$ws.Range("H129").Value = 1063.36
$ws.Range("J129").Value = 1177.6316
$ws.Range("L129").Value = 3532.8948
$ws.Range("N129").Value = -13532.8948

$ws.Range("H132").Value = 1708
$ws.Range("I132").Value = 1576
$ws.Range("K132").Value = 4728
$ws.Range("M132").Value = -2198

$ws.Range("H138").Value = 9260.875
$ws.Range("I138").Value = 1972.5625
$ws.Range("J138").Value = 12905.031
$ws.Range("K138").Value = 5917.6875
$ws.Range("L138").Value = 38715.093
$ws.Range("M138").Value = -777.6875
$ws.Range("N138").Value = -48995.093

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9329.15
$ws.Range("I61").Value = 5581.52
$ws.Range("K61").Value = 5581.52
$ws.Range("M61").Value = -5369.52

$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492

$ws.Range("H96").Value = 19500
$ws.Range("J96").Value = 19500
$ws.Range("L96").Value = 19500
$ws.Range("N96").Value = -24992

$ws.Range("H136").Value = 9329.15
$ws.Range("I136").Value = 5581.52
$ws.Range("K136").Value = 16744.56
$ws.Range("M136").Value = -14194.56

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 37399.344
$ws.Range("I134").Value = 3115.1052
$ws.Range("J134").Value = 102539.4
$ws.Range("K134").Value = 9345.3156
$ws.Range("L134").Value = 307618.2
$ws.Range("M134").Value = -6810.3156
$ws.Range("N134").Value = -312688.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 300
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1000

$ws.Range("H31").Value = 6602.604
$ws.Range("I31").Value = 6967.8857
$ws.Range("J31").Value = 5619.154
$ws.Range("K31").Value = 6967.8857
$ws.Range("L31").Value = 5619.154
$ws.Range("M31").Value = -6672.8857
$ws.Range("N31").Value = -6209.154

$ws.Range("H34").Value = 6602.604
$ws.Range("I34").Value = 6967.8857
$ws.Range("J34").Value = 5619.154
$ws.Range("K34").Value = 6967.8857
$ws.Range("L34").Value = 5619.154
$ws.Range("M34").Value = -6765.8857
$ws.Range("N34").Value = -6023.154

$ws.Range("H134").Value = 2685.7896
$ws.Range("I134").Value = 1893.3572
$ws.Range("J134").Value = 4904.6
$ws.Range("K134").Value = 5680.071599999999
$ws.Range("L134").Value = 14713.8
$ws.Range("M134").Value = -3145.071599999999
$ws.Range("N134").Value = -19783.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3670
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 3670
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 11010
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -11364

$ws.Range("H104").Value = 2290.9333
$ws.Range("J104").Value = 2575.6924
$ws.Range("L104").Value = 7727.0772
$ws.Range("N104").Value = -12969.0772

$ws.Range("H131").Value = 818.7742
$ws.Range("I131").Value = 632
$ws.Range("J131").Value = 895.1818
$ws.Range("K131").Value = 1896
$ws.Range("L131").Value = 2685.5454
$ws.Range("M131").Value = 3144
$ws.Range("N131").Value = -12765.5454

$ws.Range("H134").Value = 3748.5356
$ws.Range("I134").Value = 3543.2917
$ws.Range("J134").Value = 4980
$ws.Range("K134").Value = 10629.8751
$ws.Range("L134").Value = 14940
$ws.Range("M134").Value = -5559.875100000001
$ws.Range("N134").Value = -25080

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 846.6667
$ws.Range("I22").Value = 1060
$ws.Range("J22").Value = 633.3333
$ws.Range("K22").Value = 1060
$ws.Range("L22").Value = 633.3333
$ws.Range("M22").Value = -765
$ws.Range("N22").Value = -1223.3333

$ws.Range("H27").Value = 846.6667
$ws.Range("I27").Value = 1060
$ws.Range("J27").Value = 633.3333
$ws.Range("K27").Value = 1060
$ws.Range("L27").Value = 633.3333
$ws.Range("M27").Value = -953
$ws.Range("N27").Value = -847.3333

$ws.Range("H51").Value = 19999.334
$ws.Range("J51").Value = 19999.334
$ws.Range("L51").Value = 19999.334
$ws.Range("N51").Value = -20955.334

$ws.Range("H122").Value = 8111.8887
$ws.Range("J122").Value = 8750.5
$ws.Range("L122").Value = 26251.5
$ws.Range("N122").Value = -31151.5

$ws.Range("H132").Value = 6643.143
$ws.Range("J132").Value = 5874.5
$ws.Range("L132").Value = 17623.5
$ws.Range("N132").Value = -22683.5

$ws.Range("H141").Value = 46413.25
$ws.Range("J141").Value = 46413.25
$ws.Range("L141").Value = 46413.25
$ws.Range("N141").Value = -56773.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 30966.75
$ws.Range("J74").Value = 30966.75
$ws.Range("L74").Value = 30966.75
$ws.Range("N74").Value = -32838.75

$ws.Range("H77").Value = 30966.75
$ws.Range("J77").Value = 30966.75
$ws.Range("L77").Value = 92900.25
$ws.Range("N77").Value = -102260.25

$ws.Range("H81").Value = 20003320
$ws.Range("I81").Value = 1399
$ws.Range("J81").Value = 25003802
$ws.Range("K81").Value = 2798
$ws.Range("L81").Value = 50007604
$ws.Range("M81").Value = -1737
$ws.Range("N81").Value = -50009726

$ws.Range("H84").Value = 20003320
$ws.Range("I84").Value = 1399
$ws.Range("J84").Value = 25003802
$ws.Range("K84").Value = 13990
$ws.Range("L84").Value = 250038020
$ws.Range("M84").Value = -8686
$ws.Range("N84").Value = -250048628

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H113").Value = 624.2759
$ws.Range("I113").Value = 272.3889
$ws.Range("J113").Value = 1200.091
$ws.Range("K113").Value = 817.1667
$ws.Range("L113").Value = 3600.273
$ws.Range("M113").Value = 1352.8333
$ws.Range("N113").Value = -7940.272999999999

$ws.Range("H136").Value = 5252.314
$ws.Range("I136").Value = 2288.0334
$ws.Range("J136").Value = 9487
$ws.Range("K136").Value = 6864.100199999999
$ws.Range("L136").Value = 28461
$ws.Range("M136").Value = -4314.100199999999
$ws.Range("N136").Value = -33561

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

